$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - sheet1
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1581
$ws1.Range("F5").Value = 0
$ws1.Range("F7").Value = 405
$ws1.Range("F10").Value = 457

# Sheet "全部类型" (All types) - sheet4
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 374
$ws4.Range("F4").Value = 0
$ws4.Range("F6").Value = 0
$ws4.Range("F7").Value = 405
$ws4.Range("F9").Value = 0
$ws4.Range("F10").Value = 457
